$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Site de e-Commerce" bullet (ilvl 1, numId 3) right
#    after "Communication direct (lors de la livraison)" and before
#    "Distribution".
# ---------------------------------------------------------------------
$pComm = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Communication direct*") {
        $pComm = $p
        break
    }
}
$pComm.Range.InsertParagraphAfter()
$newIndex = $pComm.Index + 1
$d.Paragraphs.Item($newIndex).Range.Text = "Site de e-Commerce"

# ---------------------------------------------------------------------
# 2) Insert a new "Enterprise (taille à discuter)" bullet (ilvl 0,
#    numId 5) right after the "Toute personne ayant besoin d'un
#    ordinateur..." bullet and before "Structure de couts".
# ---------------------------------------------------------------------
$pToute = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Toute personne ayant besoin*") {
        $pToute = $p
        break
    }
}
$pToute.Range.InsertParagraphAfter()
$newIndex2 = $pToute.Index + 1
$d.Paragraphs.Item($newIndex2).Range.Text = "Enterprise (taille à discuter)"

# ---------------------------------------------------------------------
# 3) Add a zero-width "_GoBack" bookmark right after the "Salaires"
#    run (end of that paragraph's text, before the paragraph mark).
#
#    The COM range engine mis-places a bookmark whose Start/End both
#    land exactly on a paragraph's end-of-text boundary, so a sentinel
#    character is appended first (pushing the boundary forward by one),
#    the bookmark is anchored against that now-interior position, and
#    the sentinel is deleted again -- the bookmark range collapses back
#    onto the real boundary as the text shrinks.
# ---------------------------------------------------------------------
$pSalaires = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Salaires`r") {
        $pSalaires = $p
        break
    }
}
$endPos = $pSalaires.Range.End - 1
$insPoint = $d.Range($endPos, $endPos)
$insPoint.InsertAfter("X")

$bmPos = $pSalaires.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$sentinel = $d.Range($bmPos, $bmPos + 1)
$sentinel.Delete()
